# Update the Estonia MSME summary indicators with higher-precision figures.
# The source cells store these figures as text (not numbers), so we force
# a text number format before writing the value and then restore the
# cell's original style so no formatting is altered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B11" = "44.33"
    "C11" = "4.88"
    "D11" = "49.21"
    "B12" = "27.58"
    "C12" = "50.84"
    "D12" = "78.41"
    "B33" = "36.63"
    "C33" = "4.26"
    "D33" = "40.89"
    "B34" = "29.06"
    "C34" = "48.99"
    "D34" = "78.05"
    "B36" = "89.34"
    "D36" = "99.73"
}

foreach ($addr in $updates.Keys) {
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $updates[$addr]
    $rng.Style = $origStyle
}
